$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new time-log entry as row 14 (mirrors the existing rows 7-13) ---

# Write the new row's values first (formats are copied afterwards so the
# dependency graph for E5 = SUM(E7:E14)/60 picks up the new precedent cleanly).
$ws.Cells.Item(14, 1).Value = 41924                  # Date column (A)
$ws.Cells.Item(14, 2).Value = 0.569444444444444       # Start time (B)
$ws.Cells.Item(14, 3).Value = 0.689583333333333       # Stop time (C)
$ws.Cells.Item(14, 4).Value = 20                      # Interruption time, minutes (D)
$ws.Cells.Item(14, 5).Formula = "=((HOUR(C14)-HOUR(B14))*60)+(MINUTE(C14)-MINUTE(B14))-D14"  # Delta time (E)
$ws.Cells.Item(14, 8).Value = "Actualizaciones al plan general."                              # Comment (H)

# Copy the formatting of row 13 onto row 14 so the new row keeps the same
# per-column styles (date/time/number formats, fonts, wrap, etc.).
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)        # xlPasteFormats
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("D13:E13").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("H13").Copy()
$ws.Range("H14").PasteSpecial(-4122)

# NOTE: row 14 intentionally has no Phase/Task (F) or blank (G) entries
# (unlike row 13), so columns F/G are never touched above.

# --- Update the total Delta Time formula to include the new row ---
$ws.Cells.Item(5, 5).Formula = "=SUM(E7:E14)/60"

# --- Match the author's last selection when they saved the file ---
$ws.Range("E6").Select() | Out-Null
